$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 1698.138
$ws.Cells.Item(92, 9).Value = 1723.0714
$ws.Cells.Item(92, 10).Value = 1000
$ws.Cells.Item(92, 11).Value = 1723.0714
$ws.Cells.Item(92, 12).Value = 1000
$ws.Cells.Item(92, 13).Value = -475.0714
$ws.Cells.Item(92, 14).Value = -3496

$ws.Cells.Item(98, 8).Value = 1873.4615
$ws.Cells.Item(98, 9).Value = 1932.4
$ws.Cells.Item(98, 10).Value = 400
$ws.Cells.Item(98, 11).Value = 1932.4
$ws.Cells.Item(98, 12).Value = 400
$ws.Cells.Item(98, 13).Value = -434.4000000000001
$ws.Cells.Item(98, 14).Value = -3396

$ws.Cells.Item(118, 8).Value = 924.6
$ws.Cells.Item(118, 10).Value = 1450
$ws.Cells.Item(118, 12).Value = 4350
$ws.Cells.Item(118, 14).Value = -7664

$ws.Cells.Item(122, 8).Value = 1873.4615
$ws.Cells.Item(122, 9).Value = 1932.4
$ws.Cells.Item(122, 10).Value = 400
$ws.Cells.Item(122, 11).Value = 5797.200000000001
$ws.Cells.Item(122, 12).Value = 1200
$ws.Cells.Item(122, 13).Value = -3347.200000000001
$ws.Cells.Item(122, 14).Value = -6100

$ws.Cells.Item(137, 8).Value = 816.46344
$ws.Cells.Item(137, 9).Value = 738.2222
$ws.Cells.Item(137, 10).Value = 1379.8
$ws.Cells.Item(137, 11).Value = 2214.6666
$ws.Cells.Item(137, 12).Value = 4139.4
$ws.Cells.Item(137, 13).Value = 335.3334
$ws.Cells.Item(137, 14).Value = -9239.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 13106.937
$ws.Cells.Item(32, 9).Value = 9669.244000000001
$ws.Cells.Item(32, 10).Value = 74985.39999999999
$ws.Cells.Item(32, 11).Value = 9669.244000000001
$ws.Cells.Item(32, 12).Value = 74985.39999999999
$ws.Cells.Item(32, 13).Value = -9382.244000000001
$ws.Cells.Item(32, 14).Value = -75559.39999999999

$ws.Cells.Item(37, 8).Value = 8219
$ws.Cells.Item(37, 10).Value = 8219
$ws.Cells.Item(37, 12).Value = 8219
$ws.Cells.Item(37, 14).Value = -8765

$ws.Cells.Item(44, 8).Value = 22071.727
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 22071.727
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 22071.727
$ws.Cells.Item(44, 13).ClearContents()
$ws.Cells.Item(44, 14).Value = -23047.727

$ws.Cells.Item(55, 8).Value = 25066.262
$ws.Cells.Item(55, 10).Value = 25066.262
$ws.Cells.Item(55, 12).Value = 25066.262
$ws.Cells.Item(55, 14).Value = -25696.262

$ws.Cells.Item(61, 8).Value = 1952.3334
$ws.Cells.Item(61, 9).Value = 1614.2307
$ws.Cells.Item(61, 10).Value = 4150
$ws.Cells.Item(61, 11).Value = 1614.2307
$ws.Cells.Item(61, 12).Value = 4150
$ws.Cells.Item(61, 13).Value = -1402.2307
$ws.Cells.Item(61, 14).Value = -4574

$ws.Cells.Item(74, 8).Value = 1073.174
$ws.Cells.Item(74, 9).Value = 1043.5278
$ws.Cells.Item(74, 11).Value = 1043.5278
$ws.Cells.Item(74, 13).Value = -169.5278000000001

$ws.Cells.Item(77, 8).Value = 1073.174
$ws.Cells.Item(77, 9).Value = 1043.5278
$ws.Cells.Item(77, 11).Value = 5217.639
$ws.Cells.Item(77, 13).Value = -849.6390000000001

$ws.Cells.Item(80, 8).Value = 24005.555
$ws.Cells.Item(80, 10).Value = 26943.75
$ws.Cells.Item(80, 12).Value = 26943.75
$ws.Cells.Item(80, 14).Value = -28939.75

$ws.Cells.Item(83, 8).Value = 24005.555
$ws.Cells.Item(83, 10).Value = 26943.75
$ws.Cells.Item(83, 12).Value = 80831.25
$ws.Cells.Item(83, 14).Value = -90815.25

$ws.Cells.Item(97, 8).Value = 1769.9756
$ws.Cells.Item(97, 9).Value = 2395.8076
$ws.Cells.Item(97, 10).Value = 685.2
$ws.Cells.Item(97, 11).Value = 2395.8076
$ws.Cells.Item(97, 12).Value = 685.2
$ws.Cells.Item(97, 13).Value = -1899.8076
$ws.Cells.Item(97, 14).Value = -1677.2

$ws.Cells.Item(110, 8).Value = 967.25
$ws.Cells.Item(110, 9).Value = 848.4286
$ws.Cells.Item(110, 11).Value = 848.4286
$ws.Cells.Item(110, 13).Value = 1196.5714

$ws.Cells.Item(132, 8).Value = 1864.2
$ws.Cells.Item(132, 9).Value = 1477.5625
$ws.Cells.Item(132, 10).Value = 2815.923
$ws.Cells.Item(132, 11).Value = 4432.6875
$ws.Cells.Item(132, 12).Value = 8447.769
$ws.Cells.Item(132, 13).Value = -1902.6875
$ws.Cells.Item(132, 14).Value = -13507.769

$ws.Cells.Item(136, 8).Value = 1952.3334
$ws.Cells.Item(136, 9).Value = 1614.2307
$ws.Cells.Item(136, 10).Value = 4150
$ws.Cells.Item(136, 11).Value = 4842.6921
$ws.Cells.Item(136, 12).Value = 12450
$ws.Cells.Item(136, 13).Value = -2292.6921
$ws.Cells.Item(136, 14).Value = -17550

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 50009116
$ws.Cells.Item(86, 9).Value = 90921200
$ws.Cells.Item(86, 10).Value = 5454.1113
$ws.Cells.Item(86, 11).Value = 90921200
$ws.Cells.Item(86, 12).Value = 5454.1113
$ws.Cells.Item(86, 13).Value = -90920077
$ws.Cells.Item(86, 14).Value = -7700.1113

$ws.Cells.Item(89, 8).Value = 50009116
$ws.Cells.Item(89, 9).Value = 90921200
$ws.Cells.Item(89, 10).Value = 5454.1113
$ws.Cells.Item(89, 11).Value = 454606000
$ws.Cells.Item(89, 12).Value = 27270.5565
$ws.Cells.Item(89, 13).Value = -454600384
$ws.Cells.Item(89, 14).Value = -38502.5565

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1258.5
$ws.Cells.Item(16, 9).Value = 1021.1
$ws.Cells.Item(16, 10).Value = 1852
$ws.Cells.Item(16, 11).Value = 1021.1
$ws.Cells.Item(16, 12).Value = 1852
$ws.Cells.Item(16, 13).Value = -734.1
$ws.Cells.Item(16, 14).Value = -2426

$ws.Cells.Item(31, 8).Value = 39159.82
$ws.Cells.Item(31, 9).Value = 3324.5881
$ws.Cells.Item(31, 10).Value = 94541.55
$ws.Cells.Item(31, 11).Value = 3324.5881
$ws.Cells.Item(31, 12).Value = 94541.55
$ws.Cells.Item(31, 13).Value = -3029.5881
$ws.Cells.Item(31, 14).Value = -95131.55

$ws.Cells.Item(34, 8).Value = 39159.82
$ws.Cells.Item(34, 9).Value = 3324.5881
$ws.Cells.Item(34, 10).Value = 94541.55
$ws.Cells.Item(34, 11).Value = 3324.5881
$ws.Cells.Item(34, 12).Value = 94541.55
$ws.Cells.Item(34, 13).Value = -3122.5881
$ws.Cells.Item(34, 14).Value = -94945.55

$ws.Cells.Item(113, 8).Value = 1258.5
$ws.Cells.Item(113, 9).Value = 1021.1
$ws.Cells.Item(113, 10).Value = 1852
$ws.Cells.Item(113, 11).Value = 1021.1
$ws.Cells.Item(113, 12).Value = 1852
$ws.Cells.Item(113, 13).Value = 1148.9
$ws.Cells.Item(113, 14).Value = -6192

$ws.Cells.Item(134, 8).Value = 14286831
$ws.Cells.Item(134, 9).Value = 1164.2069
$ws.Cells.Item(134, 10).Value = 83334216
$ws.Cells.Item(134, 11).Value = 3492.620699999999
$ws.Cells.Item(134, 12).Value = 250002648
$ws.Cells.Item(134, 13).Value = -957.6206999999995
$ws.Cells.Item(134, 14).Value = -250007718

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2036.3846
$ws.Cells.Item(5, 9).Value = 1212.8235
$ws.Cells.Item(5, 10).Value = 2436.4
$ws.Cells.Item(5, 11).Value = 3638.4705
$ws.Cells.Item(5, 12).Value = 7309.200000000001
$ws.Cells.Item(5, 13).Value = -3526.4705
$ws.Cells.Item(5, 14).Value = -7533.200000000001

$ws.Cells.Item(117, 8).Value = 1962.9333
$ws.Cells.Item(117, 9).Value = 348.2
$ws.Cells.Item(117, 11).Value = 1044.6
$ws.Cells.Item(117, 13).Value = 2397.4

$ws.Cells.Item(122, 8).Value = 23971.75
$ws.Cells.Item(122, 9).Value = 565.8
$ws.Cells.Item(122, 10).Value = 26972.514
$ws.Cells.Item(122, 11).Value = 5092.2
$ws.Cells.Item(122, 12).Value = 242752.626
$ws.Cells.Item(122, 13).Value = -2642.2
$ws.Cells.Item(122, 14).Value = -247652.626

$ws.Cells.Item(135, 8).Value = 2036.3846
$ws.Cells.Item(135, 9).Value = 1212.8235
$ws.Cells.Item(135, 10).Value = 2436.4
$ws.Cells.Item(135, 11).Value = 10915.4115
$ws.Cells.Item(135, 12).Value = 21927.6
$ws.Cells.Item(135, 13).Value = -8380.4115
$ws.Cells.Item(135, 14).Value = -26997.6

$ws.Cells.Item(139, 8).Value = 35980.723
$ws.Cells.Item(139, 9).Value = 42810.875
$ws.Cells.Item(139, 11).Value = 128432.625
$ws.Cells.Item(139, 13).Value = -123292.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1645561.6
$ws.Cells.Item(122, 9).Value = 1645561.6
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4936684.800000001
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -4934234.800000001
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 4191.7144
$ws.Cells.Item(126, 10).Value = 6037.8
$ws.Cells.Item(126, 12).Value = 18113.4
$ws.Cells.Item(126, 14).Value = -23053.4

$ws.Cells.Item(134, 8).Value = 17008.834
$ws.Cells.Item(134, 10).Value = 17008.834
$ws.Cells.Item(134, 12).Value = 51026.50199999999
$ws.Cells.Item(134, 14).Value = -56096.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2484.1562
$ws.Cells.Item(68, 9).Value = 2406.15
$ws.Cells.Item(68, 10).Value = 2614.1667
$ws.Cells.Item(68, 11).Value = 2406.15
$ws.Cells.Item(68, 12).Value = 2614.1667
$ws.Cells.Item(68, 13).Value = -1657.15
$ws.Cells.Item(68, 14).Value = -4112.1667

$ws.Cells.Item(71, 8).Value = 2484.1562
$ws.Cells.Item(71, 9).Value = 2406.15
$ws.Cells.Item(71, 10).Value = 2614.1667
$ws.Cells.Item(71, 11).Value = 12030.75
$ws.Cells.Item(71, 12).Value = 13070.8335
$ws.Cells.Item(71, 13).Value = -8286.75
$ws.Cells.Item(71, 14).Value = -20558.8335

$ws.Cells.Item(107, 8).Value = 5400
$ws.Cells.Item(107, 9).Value = 5400
$ws.Cells.Item(107, 11).Value = 5400
$ws.Cells.Item(107, 13).Value = -3480

$ws.Cells.Item(132, 8).Value = 3376.6445
$ws.Cells.Item(132, 9).Value = 2504.9033
$ws.Cells.Item(132, 10).Value = 5306.9287
$ws.Cells.Item(132, 11).Value = 7514.7099
$ws.Cells.Item(132, 12).Value = 15920.7861
$ws.Cells.Item(132, 13).Value = -4984.7099
$ws.Cells.Item(132, 14).Value = -20980.7861

$ws.Cells.Item(135, 8).Value = 52379
$ws.Cells.Item(135, 10).Value = 52379
$ws.Cells.Item(135, 12).Value = 52379
$ws.Cells.Item(135, 14).Value = -62519

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 6843.727
$ws.Cells.Item(62, 9).Value = 6149.5
$ws.Cells.Item(62, 10).Value = 7240.4287
$ws.Cells.Item(62, 11).Value = 6149.5
$ws.Cells.Item(62, 12).Value = 7240.4287
$ws.Cells.Item(62, 13).Value = -5525.5
$ws.Cells.Item(62, 14).Value = -8488.4287

$ws.Cells.Item(65, 8).Value = 6843.727
$ws.Cells.Item(65, 9).Value = 6149.5
$ws.Cells.Item(65, 10).Value = 7240.4287
$ws.Cells.Item(65, 11).Value = 30747.5
$ws.Cells.Item(65, 12).Value = 36202.14350000001
$ws.Cells.Item(65, 13).Value = -27627.5
$ws.Cells.Item(65, 14).Value = -42442.14350000001

$ws.Cells.Item(96, 8).Value = 2499.7693
$ws.Cells.Item(96, 9).Value = 2317.1667
$ws.Cells.Item(96, 10).Value = 2656.2856
$ws.Cells.Item(96, 11).Value = 2317.1667
$ws.Cells.Item(96, 12).Value = 2656.2856
$ws.Cells.Item(96, 13).Value = -944.1667000000002
$ws.Cells.Item(96, 14).Value = -5402.2856

$ws.Cells.Item(122, 8).Value = 76924040
$ws.Cells.Item(122, 9).Value = 100000910
$ws.Cells.Item(122, 10).Value = 1133.3334
$ws.Cells.Item(122, 11).Value = 300002730
$ws.Cells.Item(122, 12).Value = 3400.0002
$ws.Cells.Item(122, 13).Value = -300000280
$ws.Cells.Item(122, 14).Value = -8300.0002

$ws.Cells.Item(126, 8).Value = 887.375
$ws.Cells.Item(126, 9).Value = 842.7143
$ws.Cells.Item(126, 11).Value = 2528.1429
$ws.Cells.Item(126, 13).Value = -58.14289999999983

$ws.Cells.Item(132, 8).Value = 1494.4286
$ws.Cells.Item(132, 9).Value = 1194.1482
$ws.Cells.Item(132, 10).Value = 2034.9333
$ws.Cells.Item(132, 11).Value = 3582.4446
$ws.Cells.Item(132, 12).Value = 6104.7999
$ws.Cells.Item(132, 13).Value = -1052.4446
$ws.Cells.Item(132, 14).Value = -11164.7999
